# Swap the data of row 18 and row 19 (full record swap), per the diff.
# Only columns whose values actually differ between the two rows are touched;
# columns that are identical (or empty) in both rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current ("before") values for the columns that differ ---
# (use .Value2 - in this engine .Value does not resolve to the underlying data)
$A18 = $ws.Range("A18").Value2
$B18 = $ws.Range("B18").Value2
$D18 = $ws.Range("D18").Value2
$E18 = $ws.Range("E18").Value2
$F18 = $ws.Range("F18").Value2
$G18 = $ws.Range("G18").Value2
$H18 = $ws.Range("H18").Value2
$P18 = $ws.Range("P18").Value2
$Q18 = $ws.Range("Q18").Value2
$R18 = $ws.Range("R18").Value2

$A19 = $ws.Range("A19").Value2
$B19 = $ws.Range("B19").Value2
$D19 = $ws.Range("D19").Value2
$E19 = $ws.Range("E19").Value2
$F19 = $ws.Range("F19").Value2
$G19 = $ws.Range("G19").Value2
$H19 = $ws.Range("H19").Value2
$P19 = $ws.Range("P19").Value2
$Q19 = $ws.Range("Q19").Value2
$R19 = $ws.Range("R19").Value2

# --- write row 18 with row 19's former values ---
$ws.Range("A18").Value = $A19
$ws.Range("B18").Value = $B19
$ws.Range("D18").Value = $D19
$ws.Range("E18").Value = $E19
$ws.Range("F18").Value = $F19
$ws.Range("G18").Value = $G19
$ws.Range("H18").Value = $H19
$ws.Range("J18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = "färsk spillning"
$ws.Range("P18").Value = $P19
$ws.Range("Q18").Value = $Q19
$ws.Range("R18").Value = $R19
$ws.Range("AC18").Value = ""
$ws.Range("AF18").Value = ""

# --- write row 19 with row 18's former values ---
$ws.Range("A19").Value = $A18
$ws.Range("B19").Value = $B18
$ws.Range("D19").Value = $D18
$ws.Range("E19").Value = $E18
$ws.Range("F19").Value = $F18
$ws.Range("G19").Value = $G18
$ws.Range("H19").Value = $H18
$ws.Range("J19").Value = ""
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = ""
$ws.Range("P19").Value = $P18
$ws.Range("Q19").Value = $Q18
$ws.Range("R19").Value = $R18
$ws.Range("AC19").Value = "På äldre tall."
$ws.Range("AF19").Value = ""
